$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-10 Monday", "2024-06-11 Tuesday"),
    @("738×5=", "197×9="),
    @("518×6=", "948×7="),
    @("461×3=", "227×2="),
    @("501×2=", "371×7="),
    @("179×4=", "917×4="),
    @("405×4=", "915×6="),
    @("615×3=", "230×8="),
    @("266×7=", "682×8="),
    @("922×7=", "383×5="),
    @("522×6=", "232×4="),
    @("527×6=", "264×8="),
    @("171×5=", "451×8="),
    @("269×3=", "342×6="),
    @("140×4=", "857×2="),
    @("340×2=", "484×6="),
    @("965×8=", "405×2="),
    @("233×2=", "255×3="),
    @("950×9=", "451×3="),
    @("904×3=", "272×6="),
    @("445×5=", "223×9="),
    @("285×2=", "313×9="),
    @("976×9=", "974×2="),
    @("184×8=", "624×7="),
    @("926×9=", "540×9="),
    @("652×9=", "933×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
